$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Bump the LiveSLR build number baked into the shared string for B2
$ws.Range("B2").Value = "Copyright @ 2023 Cytel Inc. LiveSLR 4.1.0.0 - Build #68318"

# Move the saved cursor/selection on Sheet1 from M7 to B2
$ws.Range("B2").Select()
